$p = $ppt.ActivePresentation

$newStyleId = "{8BF72472-F650-4209-89D6-C32640348BBC}"

foreach ($slideIdx in 14,15,16) {
    $s = $p.Slides.Item($slideIdx)
    $shape = $s.Shapes.Item(1)
    if ($shape.HasTable) {
        $shape.Table.ApplyStyle($newStyleId)
    }
}
